# Merge to new branch: add "Sheet2" with a single cell and make it the
# active/selected sheet (mirrors the target workbook.xml / sheet1.xml /
# sheet2.xml / sharedStrings.xml diff).

$wb = $excel.ActiveWorkbook

# Add a new worksheet positioned AFTER the last existing sheet (Sheet1),
# so it lands as the second tab, matching the diff's sheet order.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Sheet2"

# New shared-string value "new" written into A1 of the new sheet.
$newSheet.Range("A1").Value = "new"

# Make Sheet2 the active sheet/tab and select C4 there, which also clears
# the tabSelected flag on Sheet1's sheetView.
$newSheet.Activate()
$newSheet.Range("C4").Select()
